# Updates crypto price/volume table cells to match the latest scrape.
# A leading apostrophe forces each assignment to be stored as literal
# text (matching the source inlineStr cells) instead of being
# auto-coerced to a number by Excel's usual type inference.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.658.96"
$ws.Range("E2").Value = "'  +1.51%  "
$ws.Range("D3").Value = "'3.265.11"
$ws.Range("E3").Value = "'  +5.36%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'598.75"
$ws.Range("E5").Value = "'  +1.63%  "
$ws.Range("D6").Value = "'142.45"
$ws.Range("E6").Value = "'  +5.33%  "
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("D8").Value = "'3.265.33"
$ws.Range("E8").Value = "'  +5.58%  "
$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "'  +1.26%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "'  +2.86%  "
$ws.Range("E11").Value = "'  +2.89%  "
$ws.Range("D12").Value = "'0.470"
$ws.Range("E12").Value = "'  +3.94%  "
$ws.Range("E13").Value = "'  +1.40%  "
$ws.Range("D14").Value = "'34.79"
$ws.Range("E14").Value = "'  +3.93%  "
$ws.Range("D15").Value = "'3.800.50"
$ws.Range("E15").Value = "'  +5.42%  "
$ws.Range("E16").Value = "'  +0.79%  "
$ws.Range("D17").Value = "'3.262.49"
$ws.Range("E17").Value = "'  +5.50%  "
$ws.Range("D18").Value = "'63.700.67"
$ws.Range("E18").Value = "'  +1.53%  "
$ws.Range("D19").Value = "'6.87"
$ws.Range("E19").Value = "'  +3.85%  "
$ws.Range("D20").Value = "'479.51"
$ws.Range("E20").Value = "'  +1.56%  "
$ws.Range("D21").Value = "'14.29"
$ws.Range("E21").Value = "'  +1.09%  "
$ws.Range("D22").Value = "'0.739"
$ws.Range("E22").Value = "'  +7.29%  "
$ws.Range("D23").Value = "'8.05"
$ws.Range("E23").Value = "'  +6.52%  "
$ws.Range("D24").Value = "'13.54"
$ws.Range("E24").Value = "'  +5.47%  "
$ws.Range("D25").Value = "'84.13"
$ws.Range("E25").Value = "'  -2.55%  "
$ws.Range("E26").Value = "'  +0.01%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "'  +3.05%  "
$ws.Range("B28").Value = "'NEARProtocol"
$ws.Range("C28").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").Value = "'7.34"
$ws.Range("E28").Value = "'  +6.35%  "
$ws.Range("B29").Value = "'FirstDigitalUSD"
$ws.Range("C29").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = "'  +2.87%  "
$ws.Range("D31").Value = "'2.17"
$ws.Range("E31").Value = "'  +7.78%  "
$ws.Range("D32").Value = "'27.85"
$ws.Range("E32").Value = "'  +3.77%  "
$ws.Range("E33").Value = "'  +0.51%  "
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "'  +1.53%  "
$ws.Range("D35").Value = "'1.09"
$ws.Range("E35").Value = "'  +2.15%  "
$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = "'  +2.79%  "
$ws.Range("D37").Value = "'53.17"
$ws.Range("E37").Value = "'  +2.60%  "
$ws.Range("D38").Value = "'0.0₃0732"
$ws.Range("E38").Value = "'  +2.66%  "
$ws.Range("D39").Value = "'0.0396"
$ws.Range("E39").Value = "'  +3.28%  "
$ws.Range("D40").Value = "'421.44"
$ws.Range("E40").Value = "'  +0.90%  "
$ws.Range("B41").Value = "'dogwifhat"
$ws.Range("C41").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.79"
$ws.Range("E41").Value = "'  +3.16%  "
$ws.Range("D42").Value = "'8.44"
$ws.Range("E42").Value = "'  +3.33%  "
$ws.Range("B43").Value = "'Maker"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.990.85"
$ws.Range("E43").Value = "'  +4.99%  "
$ws.Range("E44").Value = "'  -4.46%  "
$ws.Range("D45").Value = "'0.267"
$ws.Range("E45").Value = "'  +5.45%  "
$ws.Range("D46").Value = "'2.23"
$ws.Range("E46").Value = "'  +6.39%  "
$ws.Range("B47").Value = "'InjectiveProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'26.13"
$ws.Range("E47").Value = "'  +3.85%  "
$ws.Range("B48").Value = "'ThetaToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.35"
$ws.Range("E48").Value = "'  +2.73%  "
$ws.Range("B49").Value = "'USDe"
$ws.Range("C49").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "'  +0.06%  "
$ws.Range("D50").Value = "'0.115"
$ws.Range("E50").Value = "'  +2.07%  "
$ws.Range("D51").Value = "'122.58"
$ws.Range("E51").Value = "'  +2.39%  "
